# "Now storage works and production values are not huge"
#
# The solar/wind "capacity" column (E) on Sheet1 held production figures
# that were roughly double what they should have been once storage
# round-trip losses were accounted for correctly. Rescale the affected
# rows (solar plants in rows 4-31, wind plants in rows 37-73) down to
# their corrected values, then leave the selection where the edit ended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E4").Value  = 135.970057925482
$ws.Range("E7").Value  = 116.68701050323831
$ws.Range("E10").Value = 108.77606503685921
$ws.Range("E13").Value = 135.97008129607386
$ws.Range("E16").Value = 123.60916481461271
$ws.Range("E19").Value = 151.79205439234426
$ws.Range("E22").Value = 154.26423768863671
$ws.Range("E25").Value = 140.91444788865843
$ws.Range("E28").Value = 170.58064744416549
$ws.Range("E31").Value = 182.94156392562687

$ws.Range("E37").Value = 123.95527047609363
$ws.Range("E40").Value = 155.45088567085691
$ws.Range("E43").Value = 184.32598657155032
$ws.Range("E46").Value = 129.78962305534324
$ws.Range("E49").Value = 193.42362110190592
$ws.Range("E52").Value = 126.08134811090491
$ws.Range("E55").Value = 133.00346134052322
$ws.Range("E58").Value = 154.01701935900746
$ws.Range("E61").Value = 126.2791227746083
$ws.Range("E64").Value = 235.15407514331923
$ws.Range("E67").Value = 84.993661726527705
$ws.Range("E70").Value = 158.46694929233325
$ws.Range("E73").Value = 233.17632850628524

# Reflect where the editor ended up: scrolled a little further down the
# sheet and with E73 (the last corrected cell) selected.
$ws.Activate()
$ws.Range("E73").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
